$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the cells remain formatted/stored as text (matching the original
# inline-string cell type) instead of being auto-converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "315.48"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "5"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.93%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "5"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.41%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "5"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08176"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.48%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "5"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.994"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.46%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "5"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.376"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.37%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "5"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.297"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.78%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "5"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9384"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.94%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "5"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1303"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-8.75%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "5"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1970"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.37%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "5"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09011"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.35%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "5"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03552"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.06%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "5"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09748"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.52%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "5"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001421"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.85%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "5"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006672"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "14.21%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "5"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-7.35%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "5"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.119"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-6.12%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "5"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3472"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.23%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "5"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1317"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.31%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "5"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.970"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "7.32%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "5"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2490"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.61%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "5"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04362"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.35%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "5"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001241"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.00%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "5"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004760"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "8.83%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "5"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "199.39%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "5"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-7.55%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "5"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "5"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "5"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "5"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "5"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "5"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "5"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "5"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "5"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "5"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "5"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "5"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02233"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "9.00%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "5"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05200"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.49%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "5"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "4.85%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "5"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01033"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.54%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "5"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1401"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.42%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "5"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.36%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "5"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008883"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.27%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "5"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006825"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "7.43%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "5"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.07%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "5"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003005"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "10.12%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "5"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001691"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "30.06%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "5"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.07%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "5"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.07%"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "5"
